$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value for C18 (time value of 6 hours = 0.25 of a day)
$ws.Range("C18").Value = 0.25

# Update the selected cell/range to C19 (matching the saved view state)
$ws.Range("C19").Select()
